$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

$refFont = $ws.Range("C2").Font
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("B3").ClearFormats()
$ws.Range("B3").Font.Name = $refFont.Name
$ws.Range("B3").Font.Size = $refFont.Size
$ws.Range("B3").WrapText = $true
$ws.Range("B3").VerticalAlignment = -4108

$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 18.07.2025"

$ws.Range("B6").Value = "20.07."
$ws.Range("C6").Value = "21.07."
$ws.Range("D6").Value = "KARTENZ./20.07 EDEKA RO"
$ws.Range("E6").Value = "75,63-"

$ws.Range("B7").Value = "22.07."
$ws.Range("C7").Value = "23.07."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-42224637"
$ws.Range("E7").Value = "53,65-"

$ws.Range("B8").Value = "23.07."
$ws.Range("C8").Value = "24.07."
$ws.Range("D8").Value = "EBAY MKTPLC EU HWDMKS"
$ws.Range("E8").Value = "178,36-"

$ws.Range("B9").Value = "24.07."
$ws.Range("C9").Value = "25.07."
$ws.Range("D9").Value = "BURGER KING Badoberan"
$ws.Range("E9").Value = "15,76-"

$ws.Range("B10").Value = "27.07."
$ws.Range("C10").Value = "28.07."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E10").Value = "24,85-"
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4107
$ws.Range("E10").WrapText = $false

$ws.Range("D12").Value = "KONTOSTAND AM 31.07.2025"
$ws.Range("E12").Value = "348,25-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 06.08.2025"
